# Excess mortality provinces - Week 48 (2022) update
# - Minor upward revisions to several previously reported weeks (rows 129-152)
# - Two newly published weeks appended: 2022 week 47 (row 153) and
#   2022 week 48 (row 154), including their "percentage change vs 2020"
#   helper formulas in columns AE:AP
# - View state: scrolled down / new active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("excess_mortality_provinces")
$ws.Activate()

# --- Revisions to previously reported weeks (129-152) ---
$ws.Range("T129").Value = 60
$ws.Range("X129").Value = 594
$ws.Range("W130").Value = 450
$ws.Range("X130").Value = 611
$ws.Range("X131").Value = 709
$ws.Range("W137").Value = 456
$ws.Range("W140").Value = 445
$ws.Range("AA140").Value = 232
$ws.Range("X141").Value = 600
$ws.Range("Q142").Value = 123
$ws.Range("X142").Value = 596
$ws.Range("U143").Value = 375
$ws.Range("Z143").Value = 419
$ws.Range("X144").Value = 609
$ws.Range("X145").Value = 661
$ws.Range("U146").Value = 448
$ws.Range("U147").Value = 460
$ws.Range("W147").Value = 478
$ws.Range("X147").Value = 692
$ws.Range("Z147").Value = 487
$ws.Range("AA147").Value = 258
$ws.Range("U148").Value = 420
$ws.Range("V148").Value = 239
$ws.Range("X148").Value = 663
$ws.Range("S149").Value = 225
$ws.Range("Z149").Value = 500
$ws.Range("U150").Value = 391
$ws.Range("V150").Value = 235
$ws.Range("W150").Value = 479
$ws.Range("X150").Value = 660
$ws.Range("Z150").Value = 462
$ws.Range("AA150").Value = 242
$ws.Range("Q151").Value = 138
$ws.Range("S151").Value = 216
$ws.Range("T151").Value = 63
$ws.Range("U151").Value = 417
$ws.Range("V151").Value = 236
$ws.Range("W151").Value = 441
$ws.Range("X151").Value = 612
$ws.Range("Y151").Value = 74
$ws.Range("Z151").Value = 477
$ws.Range("AA151").Value = 246
$ws.Range("P152").Value = 109
$ws.Range("Q152").Value = 141
$ws.Range("R152").Value = 93
$ws.Range("S152").Value = 222
$ws.Range("T152").Value = 58
$ws.Range("U152").Value = 408
$ws.Range("V152").Value = 232
$ws.Range("W152").Value = 467
$ws.Range("X152").Value = 649
$ws.Range("Y152").Value = 79
$ws.Range("Z152").Value = 486
$ws.Range("AA152").Value = 245

# --- New data: row 153 = 2022 week 47 ---
$ws.Range("N153").Value = 2022
$ws.Range("O153").Value = 47
$ws.Range("P153").Value = 123
$ws.Range("Q153").Value = 141
$ws.Range("R153").Value = 123
$ws.Range("S153").Value = 219
$ws.Range("T153").Value = 61
$ws.Range("U153").Value = 429
$ws.Range("V153").Value = 249
$ws.Range("W153").Value = 484
$ws.Range("X153").Value = 669
$ws.Range("Y153").Value = 104
$ws.Range("Z153").Value = 496
$ws.Range("AA153").Value = 267
$ws.Range("AC153").Value = 2022
$ws.Range("AD153").Value = 47
$ws.Range("AE153").Formula = "=ROUND((P153-B153)/B153*100,2)"
$ws.Range("AF153").Formula = "=ROUND((Q153-C153)/C153*100,2)"
$ws.Range("AG153").Formula = "=ROUND((R153-D153)/D153*100,2)"
$ws.Range("AH153").Formula = "=ROUND((S153-E153)/E153*100,2)"
$ws.Range("AI153").Formula = "=ROUND((T153-F153)/F153*100,2)"
$ws.Range("AJ153").Formula = "=ROUND((U153-G153)/G153*100,2)"
$ws.Range("AK153").Formula = "=ROUND((V153-H153)/H153*100,2)"
$ws.Range("AL153").Formula = "=ROUND((W153-I153)/I153*100,2)"
$ws.Range("AM153").Formula = "=ROUND((X153-J153)/J153*100,2)"
$ws.Range("AN153").Formula = "=ROUND((Y153-K153)/K153*100,2)"
$ws.Range("AO153").Formula = "=ROUND((Z153-L153)/L153*100,2)"
$ws.Range("AP153").Formula = "=ROUND((AA153-M153)/M153*100,2)"

# --- New data: row 154 = 2022 week 48 ---
$ws.Range("N154").Value = 2022
$ws.Range("O154").Value = 48
$ws.Range("P154").Value = 128
$ws.Range("Q154").Value = 154
$ws.Range("R154").Value = 121
$ws.Range("S154").Value = 216
$ws.Range("T154").Value = 62
$ws.Range("U154").Value = 406
$ws.Range("V154").Value = 222
$ws.Range("W154").Value = 534
$ws.Range("X154").Value = 700
$ws.Range("Y154").Value = 75
$ws.Range("Z154").Value = 535
$ws.Range("AA154").Value = 251
$ws.Range("AC154").Value = 2022
$ws.Range("AD154").Value = 48
$ws.Range("AE154").Formula = "=ROUND((P154-B154)/B154*100,2)"
$ws.Range("AF154").Formula = "=ROUND((Q154-C154)/C154*100,2)"
$ws.Range("AG154").Formula = "=ROUND((R154-D154)/D154*100,2)"
$ws.Range("AH154").Formula = "=ROUND((S154-E154)/E154*100,2)"
$ws.Range("AI154").Formula = "=ROUND((T154-F154)/F154*100,2)"
$ws.Range("AJ154").Formula = "=ROUND((U154-G154)/G154*100,2)"
$ws.Range("AK154").Formula = "=ROUND((V154-H154)/H154*100,2)"
$ws.Range("AL154").Formula = "=ROUND((W154-I154)/I154*100,2)"
$ws.Range("AM154").Formula = "=ROUND((X154-J154)/J154*100,2)"
$ws.Range("AN154").Formula = "=ROUND((Y154-K154)/K154*100,2)"
$ws.Range("AO154").Formula = "=ROUND((Z154-L154)/L154*100,2)"
$ws.Range("AP154").Formula = "=ROUND((AA154-M154)/M154*100,2)"

# --- View state: scroll position + active selection ---
$excel.ActiveWindow.ScrollRow = 95
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("AI155").Select()
